$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Forming the consolidated report: fill in the computed "Absent" (H) values
# for rows that previously had no value / stale values.
$ws.Range("H3").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 0
